$wb = $excel.ActiveWorkbook

# ----- Sheet 1: Overview -----
$ws = $wb.Worksheets.Item(1)

# Set cell values (drives sharedStrings content)
$ws.Range('A1').Value = 'File Name'
$ws.Range('B1').Value = 'zh-cn'
$ws.Range('C1').Value = 'de-de'
$ws.Range('A2').Value = 'ffffb80e2993-1aae-4258-b44b-7f85ee356543.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = 'Handed back: in sync with en-US'
$ws.Range('A3').Value = 'ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = 'Handed back: in sync with en-US'
$ws.Range('A4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.md'
$ws.Range('B4').Value = 'Ready for handoff'
$ws.Range('C4').Value = 'Ready for handoff'
$ws.Range('A5').Value = '.localization-config'
$ws.Range('B5').Value = 'Not to be localized'
$ws.Range('C5').Value = 'Not to be localized'

# Rebuild hyperlinks in the same ref/target order as before, with updated display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/81ddff20-c635-4825-9f13-3fbb5b8226c2.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffb80e2993-1aae-4258-b44b-7f85ee356543.md')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/ffffb80e2993-1aae-4258-b44b-7f85ee356543.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.md')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/.localization-config', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.localization-config')

# ----- Sheet 2: zh-cn -----
$ws = $wb.Worksheets.Item(2)

# Set cell values (drives sharedStrings content)
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'Status'
$ws.Range('C1').Value = 'Latest Handoff File'
$ws.Range('D1').Value = 'Latest Handoff Datetime'
$ws.Range('E1').Value = 'Latest Target File'
$ws.Range('F1').Value = 'Latest Handback File'
$ws.Range('G1').Value = 'Latest Handback DateTime'
$ws.Range('H1').Value = 'Handoff Reason'
$ws.Range('I1').Value = 'Dependency From'
$ws.Range('A2').Value = 'ffffb80e2993-1aae-4258-b44b-7f85ee356543.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf'
$ws.Range('D2').Value = '2016-02-22 17:56:45'
$ws.Range('E2').Value = 'f13d269d-2389-4baf-9322-ab170051d945.md'
$ws.Range('F2').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf'
$ws.Range('G2').Value = '2016-02-22 17:57:26'
$ws.Range('H2').Value = 'Include'
$ws.Range('A3').Value = 'ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf'
$ws.Range('D3').Value = '2016-02-22 17:56:45'
$ws.Range('E3').Value = 'f13d269d-2389-4baf-9322-ab170051d945.md'
$ws.Range('F3').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf'
$ws.Range('G3').Value = '2016-02-22 17:57:26'
$ws.Range('H3').Value = 'Include'
$ws.Range('A4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.md'
$ws.Range('B4').Value = 'Ready for handoff'
$ws.Range('C4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf'
$ws.Range('D4').Value = '2016-02-22 18:02:51'
$ws.Range('E4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.md'
$ws.Range('F4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf'
$ws.Range('G4').Value = '2016-02-22 18:01:57'
$ws.Range('H4').Value = 'Include'
$ws.Range('A5').Value = '.localization-config'
$ws.Range('B5').Value = 'Not to be localized'
$ws.Range('D5').Value = '0001-01-01 00:00:00'
$ws.Range('G5').Value = '0001-01-01 00:00:00'
$ws.Range('H5').Value = 'Ignored'

# Rebuild hyperlinks in the same ref/target order as before, with updated display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/81ddff20-c635-4825-9f13-3fbb5b8226c2.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffb80e2993-1aae-4258-b44b-7f85ee356543.md')
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f855b00d49a1a181671a52b83f46c69feeadc4b6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b47cbdc0fa9b79e31bf1a41c6fc3dd6339cd3130/e2e/81ddff20-c635-4825-9f13-3fbb5b8226c2.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.md')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/0ab9a82557735ada995cddfd047ddff72a93bb65/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/ffffb80e2993-1aae-4258-b44b-7f85ee356543.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec4d7c626e85b0a8e1d3603047f57a38cc8b313b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('E3'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fc543401a66968563956a0f264f6d20e44acd621/e2e/f13d269d-2389-4baf-9322-ab170051d945.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.md')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/05ea5889f81f2b28fcdd391da9dc00e6225390e6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.md')
$ws.Hyperlinks.Add($ws.Range('C4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ec4d7c626e85b0a8e1d3603047f57a38cc8b313b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('E4'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fc543401a66968563956a0f264f6d20e44acd621/e2e/f13d269d-2389-4baf-9322-ab170051d945.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.md')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/05ea5889f81f2b28fcdd391da9dc00e6225390e6/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/.localization-config', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.localization-config')

# ----- Sheet 3: de-de -----
$ws = $wb.Worksheets.Item(3)

# Set cell values (drives sharedStrings content)
$ws.Range('A1').Value = 'Source File Name'
$ws.Range('B1').Value = 'Status'
$ws.Range('C1').Value = 'Latest Handoff File'
$ws.Range('D1').Value = 'Latest Handoff Datetime'
$ws.Range('E1').Value = 'Latest Target File'
$ws.Range('F1').Value = 'Latest Handback File'
$ws.Range('G1').Value = 'Latest Handback DateTime'
$ws.Range('H1').Value = 'Handoff Reason'
$ws.Range('I1').Value = 'Dependency From'
$ws.Range('A2').Value = 'ffffb80e2993-1aae-4258-b44b-7f85ee356543.md'
$ws.Range('B2').Value = 'Handed back: in sync with en-US'
$ws.Range('C2').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf'
$ws.Range('D2').Value = '2016-02-22 17:56:56'
$ws.Range('E2').Value = 'f13d269d-2389-4baf-9322-ab170051d945.md'
$ws.Range('F2').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf'
$ws.Range('G2').Value = '2016-02-22 17:57:45'
$ws.Range('H2').Value = 'Include'
$ws.Range('A3').Value = 'ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md'
$ws.Range('B3').Value = 'Handed back: in sync with en-US'
$ws.Range('C3').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf'
$ws.Range('D3').Value = '2016-02-22 17:56:56'
$ws.Range('E3').Value = 'f13d269d-2389-4baf-9322-ab170051d945.md'
$ws.Range('F3').Value = 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf'
$ws.Range('G3').Value = '2016-02-22 17:57:45'
$ws.Range('H3').Value = 'Include'
$ws.Range('A4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.md'
$ws.Range('B4').Value = 'Ready for handoff'
$ws.Range('C4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf'
$ws.Range('D4').Value = '2016-02-22 18:03:04'
$ws.Range('E4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.md'
$ws.Range('F4').Value = '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf'
$ws.Range('G4').Value = '2016-02-22 18:02:18'
$ws.Range('H4').Value = 'Include'
$ws.Range('A5').Value = '.localization-config'
$ws.Range('B5').Value = 'Not to be localized'
$ws.Range('D5').Value = '0001-01-01 00:00:00'
$ws.Range('G5').Value = '0001-01-01 00:00:00'
$ws.Range('H5').Value = 'Ignored'

# Rebuild hyperlinks in the same ref/target order as before, with updated display text
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/81ddff20-c635-4825-9f13-3fbb5b8226c2.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffb80e2993-1aae-4258-b44b-7f85ee356543.md')
$ws.Hyperlinks.Add($ws.Range('C2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3d71af7fe907604226c4ecd4d569131232a4018e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('E2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/9216a8fd405693dd3007804ee05a273a684305ee/e2e/81ddff20-c635-4825-9f13-3fbb5b8226c2.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.md')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/6a3eddd3fc106a0fb1c2511027b3c3ac62d89a08/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/ffffb80e2993-1aae-4258-b44b-7f85ee356543.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md')
$ws.Hyperlinks.Add($ws.Range('C3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e50deee38aebe23fc56cea68436aa17ad67e82f0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('E3'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d45b64731103203cc1e56a2a28b8d0ebbb5d9816/e2e/f13d269d-2389-4baf-9322-ab170051d945.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.md')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/bbac462fb60e1c63e07a6f3b35f69a472552de91/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/e2e/ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.md')
$ws.Hyperlinks.Add($ws.Range('C4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e50deee38aebe23fc56cea68436aa17ad67e82f0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('E4'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d45b64731103203cc1e56a2a28b8d0ebbb5d9816/e2e/f13d269d-2389-4baf-9322-ab170051d945.md', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.md')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/bbac462fb60e1c63e07a6f3b35f69a472552de91/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/a5d64095f8bfff513b8091f06e35ff1e4e93b7b3/.localization-config', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, '.localization-config')
